$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A180").Value = "13-09-2021"
$ws.Range("B180").Value = 3.38
$ws.Range("C180").Value = 3.32
$ws.Range("D180").Value = 3.25

$ws.Range("A181").Value = "14-09-2021"
$ws.Range("B181").Value = 3.36
$ws.Range("C181").Value = 3.26
$ws.Range("D181").Value = 3.31

$ws.Range("A182").Value = "15-09-2021"
$ws.Range("B182").Value = 3.35
$ws.Range("C182").Value = 3.18
$ws.Range("D182").Value = 3.31

$ws.Range("A183").Value = "16-09-2021"
$ws.Range("D183").Value = 3.32
